$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 127, pushing the existing
# rows 127-151 down to 129-153.
$ws.Rows.Item(127).Insert()
$ws.Rows.Item(127).Insert()

# New row 127: Primera, $/caja 7 kilos, Provincia de Diguillín
$ws.Range("A127").Value = 7
$ws.Range("B127").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C127").Value = "Ñuble"
$ws.Range("D127").Value = 44505
$ws.Range("E127").Value = 16
$ws.Range("F127").Value = "Fruta"
$ws.Range("G127").Value = 100101
$ws.Range("H127").Value = "Berries"
$ws.Range("I127").Value = 100112025
$ws.Range("J127").Value = "Frutilla"
$ws.Range("K127").Value = "Sin especificar"
$ws.Range("L127").Value = "Primera"
$ws.Range("M127").Value = 80
$ws.Range("N127").Value = 7000
$ws.Range("O127").Value = 7500
$ws.Range("P127").Value = 7250
$ws.Range("Q127").Value = "`$/caja 7 kilos"
$ws.Range("R127").Value = "Provincia de Diguillín"
$ws.Range("S127").Value = 1036
$ws.Range("T127").Value = 7

# New row 128: Segunda, $/caja 7 kilos, Provincia de Diguillín
$ws.Range("A128").Value = 7
$ws.Range("B128").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C128").Value = "Ñuble"
$ws.Range("D128").Value = 44505
$ws.Range("E128").Value = 16
$ws.Range("F128").Value = "Fruta"
$ws.Range("G128").Value = 100101
$ws.Range("H128").Value = "Berries"
$ws.Range("I128").Value = 100112025
$ws.Range("J128").Value = "Frutilla"
$ws.Range("K128").Value = "Sin especificar"
$ws.Range("L128").Value = "Segunda"
$ws.Range("M128").Value = 80
$ws.Range("N128").Value = 6000
$ws.Range("O128").Value = 6500
$ws.Range("P128").Value = 6250
$ws.Range("Q128").Value = "`$/caja 7 kilos"
$ws.Range("R128").Value = "Provincia de Diguillín"
$ws.Range("S128").Value = 893
$ws.Range("T128").Value = 7
